$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "16/12/2019"
$ws.Range("B11").Value = "Mock review was conducted. "
$ws.Range("C11").Value = "Title was choosen"
$ws.Range("C12").Value = "Changes were made to the document accordingly based on the review"

$ws.Range("D30:D31").Select()
